{"js": "// Update the PROFILE summary: \"solution-oriented Software Developer\" ->\n// \"solution-driven Software Engineer\", and add \"C#, MySQL\" to the list of\n// proficient back-end languages mentioned in the summary sentence.\nconst body = context.document.body;\n\nconst titleSearch = body.search(\"solution-oriented Software Developer\", { matchCase: true });\ntitleSearch.load(\"items\");\nawait context.sync();\nif (titleSearch.items.length > 0) {\n  titleSearch.items[0].insertText(\"solution-driven Software Engineer\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst langSearch = body.search(\"proficient in Java, C, and React,\", { matchCase: true });\nlangSearch.load(\"items\");\nawait context.sync();\nif (langSearch.items.length > 0) {\n  langSearch.items[0].insertText(\"proficient in Java, C, C#, MySQL and React,\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// TECHNICAL SKILLS / Front-End: add \"React JS\" and \"Next JS\" before \"React Native\".\nconst frontEndSearch = body.search(\"React, React Native\", { matchCase: true });\nfrontEndSearch.load(\"items\");\nawait context.sync();\nif (frontEndSearch.items.length > 0) {\n  frontEndSearch.items[0].insertText(\"React JS, Next JS React Native\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// TECHNICAL SKILLS / Back-End: append \"(Proficient with MySQL workbench)\" with\n// \"MySQL workbench\" in bold, right after the existing \", C#, MySQL\" text.\n// (Use a longer, unique anchor -- plain \", C#, MySQL\" also now occurs in the\n// PROFILE paragraph above after the edit made there.)\nconst backEndSearch = body.search(\"C++ (beginner), C#, MySQL\", { matchCase: true });\nbackEndSearch.load(\"items\");\nawait context.sync();\nif (backEndSearch.items.length > 0) {\n  backEndSearch.items[0].insertText(\n    \"C++ (beginner), C#, MySQL (Proficient with MySQL workbench)\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  const boldSearch = body.search(\"MySQL workbench\", { matchCase: true });\n  boldSearch.load(\"items\");\n  await context.sync();\n  if (boldSearch.items.length > 0) {\n    boldSearch.items[0].font.bold = true;\n    await context.sync();\n  }\n}\n", "ps1": "# Resume update: refresh the PROFILE summary wording and extend the\n# TECHNICAL SKILLS bullet lines.\n$d = $word.ActiveDocument\n\n# PROFILE: \"solution-oriented Software Developer\" -> \"solution-driven Software Engineer\"\n$find = $d.Content.Find\n$find.Text = \"solution-oriented Software Developer\"\n$find.Replacement.Text = \"solution-driven Software Engineer\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# PROFILE: mention C# and MySQL alongside Java/C in the proficiency sentence.\n$find = $d.Content.Find\n$find.Text = \"proficient in Java, C, and React,\"\n$find.Replacement.Text = \"proficient in Java, C, C#, MySQL and React,\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# TECHNICAL SKILLS / Front-End: add \"React JS\" and \"Next JS\" ahead of \"React Native\".\n$find = $d.Content.Find\n$find.Text = \"React, React Native\"\n$find.Replacement.Text = \"React JS, Next JS React Native\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# TECHNICAL SKILLS / Back-End: note proficiency with MySQL Workbench.\n$find = $d.Content.Find\n$find.Text = \"C++ (beginner), C#, MySQL\"\n$find.Replacement.Text = \"C++ (beginner), C#, MySQL (Proficient with MySQL workbench)\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# Make just the \"MySQL workbench\" phrase bold.\n$boldRange = $d.Content\n$boldFind = $boldRange.Find\n$boldFind.Text = \"MySQL workbench\"\n$boldFind.Execute()\n$boldRange.Bold = 1\n"}
